# Fruta / hortaliza, semanal
# The weekly refresh rotates the price-report rows: each of rows
# 2, 4, 5, 6, 8 takes on the values previously held by another row
# in the same set (dates, quality, volume, prices, unit, origin, etc.),
# per the supplied diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for the rows that get rotated (2, 4, 5, 6, 8)
# columns D, L, M, N, O, P, Q, R, S, T
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")
$rows = @(2, 4, 5, 6, 8)

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Target row gets source row's pre-edit data
$rowMap = @{
    2 = 8
    4 = 6
    5 = 2
    6 = 5
    8 = 4
}

foreach ($target in $rowMap.Keys) {
    $source = $rowMap[$target]
    $data = $snapshot[$source]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value = $data[$c]
    }
}
